$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-22 down to 17-23
$ws.Rows.Item(16).Insert()

# Copy static columns (A,B,C,E,F,G,H,I,N,O,Q,R) from the row now at 17 (old row 16) into new row 16
$cols = @(1,2,3,5,6,7,8,9,14,15,17,18)
foreach ($c in $cols) {
    $ws.Cells.Item(16, $c).Value2 = $ws.Cells.Item(17, $c).Value2
}

# Match the date cell format (style index 2, used by column D) from row 17 onto new row 16
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat

# Set the new row 16 values
$ws.Cells.Item(16, 4).Value = 44468
$ws.Cells.Item(16, 10).Value = 300
$ws.Cells.Item(16, 11).Value = 900
$ws.Cells.Item(16, 12).Value = 1000
$ws.Cells.Item(16, 13).Value = 950
$ws.Cells.Item(16, 16).Value = 475
